# SharePoint Administrator Checklist: add Controls 4.6 and 4.7
# (gap-analysis update) and push the trailing "FSI Agent Governance
# Framework v1.0 Beta" banner down two rows to keep blank spacer rows
# above it, matching the new dimension A1:E13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Move the "FSI Agent Governance Framework v1.0 Beta" banner from
#     row 11 down to row 13 (value + formatting), then clear the old
#     location, before we touch anything else in that range.
$ws.Range("A11").Copy($ws.Range("A13"))
$ws.Range("A11").Clear()

# --- New Control 4.6 (row 9) ---
# Leading apostrophes keep "4.6"/"4.7" as text (matching the existing
# "4.1".."4.5" control-id cells) instead of being parsed as numbers.
$ws.Range("A9").Value = "'4.6"
$ws.Range("B9").Value = "Grounding Scope Governance"
$ws.Range("C9").Value = "Not Started"

# --- New Control 4.7 (row 10) ---
$ws.Range("A10").Value = "'4.7"
$ws.Range("B10").Value = "Microsoft 365 Copilot Data Governance"
$ws.Range("C10").Value = "Not Started"

# --- Blank spacer rows: row 2 (under the title) and rows 11-12
#     (above the banner, now at row 13) ---
$ws.Rows.Item(2).OutlineLevel = 0
$ws.Rows.Item(11).OutlineLevel = 0
$ws.Rows.Item(12).OutlineLevel = 0
